$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert two new blank rows at row 4 (pushes old rows 4-29 down to
#    rows 6-31, and Excel auto-updates the used-range dimension).
# ------------------------------------------------------------------
$ws.Rows("4:5").Insert()

# Copy the A:B column formatting (bold/border/centered index style,
# plain style for the label cell) from row 2 down onto the two new
# rows so the new A/B cells match the sheet's existing look.
$ws.Range("A2:B2").Copy()
$ws.Range("A4:B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Populate the two new rows: index column, method-name label, and
#    the freshly (re)simulated statistics for each of the 18 metric
#    columns C:T.
# ------------------------------------------------------------------
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"

$row4Vals = 0.9695264157487977, 2.340520323725387, 0.8132570485170263, 0.8991765622141452, 0.8991765622141452, 0.7838669524610974, 0.7838669524610974, 1.017299194204401, 0.8991765622141452, 1.017299194204401, 0.9005830733327489, 0.9005830733327489, 0.8714743983941746, 0.9001142362932143, 0.9001142362932143, 0.8998798177734471, 0.8998798177734471, 1.137274416145142
$row5Vals = 0.1400439674751711, 4.594734806660918, 0.0354867531263421, 1.880932696052732, 1.880932696052732, 3.770219205583233, 3.770219205583233, 1.909539214169965, 1.880932696052732, 1.909539214169965, 2.839879209876599, 2.839879209876599, 1.90508172429318, 2.52023037193531, 2.52023037193531, 2.360405952964665, 2.360405952964665, 2.055159440511393

for ($i = 0; $i -lt 18; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $row4Vals[$i]
    $ws.Cells.Item(5, 3 + $i).Value = $row5Vals[$i]
}

# ------------------------------------------------------------------
# 3) Rename the "Thomas Hex" method to "Matthies Hex". After the
#    insert above, that row (originally row 9, old index 7) now lives
#    at row 11.
# ------------------------------------------------------------------
$ws.Range("B11").Value = "Matthies Hex"

Write-Output "edit applied"
